$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column E ("is_new") is TRUE; these rows only had data through
# column F ("cued") and are missing the image1/type1..image4/type4 columns
# (G..N). Fill those with "N/A" to match the rest of the working set.
$rows = @(2,9,11,16,20,27,31,34,36,38,45,49,50,54,59,62,63,65,67,68,71,72,73,78,80,81,83,84,85,87,88,89,91,93,96,98,101,103,105,108,110,111,117,124,129,133,134,138,141,145,147,151,154,158,162,164,166,169,170,173,175,181,187,192)

foreach ($r in $rows) {
    for ($c = 7; $c -le 14; $c++) {
        $ws.Cells.Item($r, $c).Value = "N/A"
    }
}
